# Update cached profit/loss figures (columns H:N) on several rows across
# the ALC, ARM, BSM, CRP, CUL, GSM and LTW sheets to reflect refreshed
# market data. A handful of cells (L/M on a few rows) no longer apply and
# are cleared outright rather than set to a value.
$wb = $excel.ActiveWorkbook

# Sheet 1: ALC
$ws = $wb.Worksheets.Item(1)
$ws.Range("H64").Value = 2982.2222
$ws.Range("I64").Value = 2823.077
$ws.Range("J64").Value = 3046.875
$ws.Range("K64").Value = 2823.077
$ws.Range("L64").Value = 3046.875
$ws.Range("M64").Value = -2575.077
$ws.Range("N64").Value = -3542.875
$ws.Range("H67").Value = 2982.2222
$ws.Range("I67").Value = 2823.077
$ws.Range("J67").Value = 3046.875
$ws.Range("K67").Value = 2823.077
$ws.Range("L67").Value = 3046.875
$ws.Range("M67").Value = -1965.077
$ws.Range("N67").Value = -4762.875
$ws.Range("H74").Value = 7159.52
$ws.Range("I74").Value = 12620
$ws.Range("J74").Value = 3519.2
$ws.Range("K74").Value = 12620
$ws.Range("L74").Value = 3519.2
$ws.Range("M74").Value = -11684
$ws.Range("N74").Value = -5391.2
$ws.Range("H76").Value = 7403.057
$ws.Range("I76").Value = 10167.056
$ws.Range("K76").Value = 10167.056
$ws.Range("M76").Value = -9852.056
$ws.Range("H77").Value = 7159.52
$ws.Range("I77").Value = 12620
$ws.Range("J77").Value = 3519.2
$ws.Range("K77").Value = 63100
$ws.Range("L77").Value = 17596
$ws.Range("M77").Value = -58420
$ws.Range("N77").Value = -26956
$ws.Range("H79").Value = 7403.057
$ws.Range("I79").Value = 10167.056
$ws.Range("K79").Value = 10167.056
$ws.Range("M79").Value = -9075.056

# Sheet 2: ARM
$ws = $wb.Worksheets.Item(2)
$ws.Range("H2").Value = 1975.2354
$ws.Range("I2").Value = 1293
$ws.Range("J2").Value = 3226
$ws.Range("K2").Value = 1293
$ws.Range("L2").Value = 3226
$ws.Range("M2").Value = -1180
$ws.Range("N2").Value = -3452
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("L97").ClearContents()
$ws.Range("M97").ClearContents()
$ws.Range("N97").Value = 0
$ws.Range("H116").Value = 1975.2354
$ws.Range("I116").Value = 1293
$ws.Range("J116").Value = 3226
$ws.Range("K116").Value = 1293
$ws.Range("L116").Value = 3226
$ws.Range("M116").Value = 1001
$ws.Range("N116").Value = -7814

# Sheet 3: BSM
$ws = $wb.Worksheets.Item(3)
$ws.Range("H3").Value = 1975.2354
$ws.Range("I3").Value = 1293
$ws.Range("J3").Value = 3226
$ws.Range("K3").Value = 1293
$ws.Range("L3").Value = 3226
$ws.Range("M3").Value = -1179
$ws.Range("N3").Value = -3454
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("L29").ClearContents()
$ws.Range("M29").ClearContents()
$ws.Range("N29").Value = 0
$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("L94").ClearContents()
$ws.Range("M94").ClearContents()
$ws.Range("N94").Value = 0
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").ClearContents()
$ws.Range("N98").Value = 0

# Sheet 4: CRP
$ws = $wb.Worksheets.Item(4)
$ws.Range("H62").Value = 3544.4
$ws.Range("I62").Value = 2611.5
$ws.Range("J62").Value = 4166.3335
$ws.Range("K62").Value = 2611.5
$ws.Range("L62").Value = 4166.3335
$ws.Range("M62").Value = -1987.5
$ws.Range("N62").Value = -5414.3335
$ws.Range("H65").Value = 3544.4
$ws.Range("I65").Value = 2611.5
$ws.Range("J65").Value = 4166.3335
$ws.Range("K65").Value = 13057.5
$ws.Range("L65").Value = 20831.6675
$ws.Range("M65").Value = -9937.5
$ws.Range("N65").Value = -27071.6675
$ws.Range("H131").Value = 40571
$ws.Range("J131").Value = 40571
$ws.Range("L131").Value = 40571
$ws.Range("N131").Value = -50651

# Sheet 5: CUL
$ws = $wb.Worksheets.Item(5)
$ws.Range("H34").Value = 1301.125
$ws.Range("I34").Value = 322
$ws.Range("J34").Value = 2933
$ws.Range("K34").Value = 966
$ws.Range("L34").Value = 8799
$ws.Range("M34").Value = -882
$ws.Range("N34").Value = -8967
$ws.Range("H104").Value = 33336796
$ws.Range("I104").Value = 707.3333
$ws.Range("J104").Value = 47623692
$ws.Range("K104").Value = 2121.9999
$ws.Range("L104").Value = 142871076
$ws.Range("M104").Value = 499.0001000000002
$ws.Range("N104").Value = -142876318
$ws.Range("H122").Value = 703.1070999999999
$ws.Range("I122").Value = 348.375
$ws.Range("J122").Value = 1176.0834
$ws.Range("K122").Value = 3135.375
$ws.Range("L122").Value = 10584.7506
$ws.Range("M122").Value = -685.375
$ws.Range("N122").Value = -15484.7506

# Sheet 6: GSM
$ws = $wb.Worksheets.Item(6)
$ws.Range("H70").Value = 30310.426
$ws.Range("I70").Value = 37787.098
$ws.Range("J70").Value = 4557.4443
$ws.Range("K70").Value = 37787.098
$ws.Range("L70").Value = 4557.4443
$ws.Range("M70").Value = -37517.098
$ws.Range("N70").Value = -5097.4443
$ws.Range("H73").Value = 30310.426
$ws.Range("I73").Value = 37787.098
$ws.Range("J73").Value = 4557.4443
$ws.Range("K73").Value = 37787.098
$ws.Range("L73").Value = 4557.4443
$ws.Range("M73").Value = -36851.098
$ws.Range("N73").Value = -6429.4443
$ws.Range("H97").Value = 1581.75
$ws.Range("I97").Value = 1581.4286
$ws.Range("J97").Value = 1582.2
$ws.Range("K97").Value = 1581.4286
$ws.Range("L97").Value = 1582.2
$ws.Range("M97").Value = -1085.4286
$ws.Range("N97").Value = -2574.2
$ws.Range("H132").Value = 1775.5416
$ws.Range("I132").Value = 1716.5588
$ws.Range("J132").Value = 1918.7858
$ws.Range("K132").Value = 5149.6764
$ws.Range("L132").Value = 5756.357400000001
$ws.Range("M132").Value = -2619.6764
$ws.Range("N132").Value = -10816.3574

# Sheet 7: LTW
$ws = $wb.Worksheets.Item(7)
$ws.Range("H68").Value = 1818.9524
$ws.Range("I68").Value = 1809.9
$ws.Range("J68").Value = 2000
$ws.Range("K68").Value = 1809.9
$ws.Range("L68").Value = 2000
$ws.Range("M68").Value = -1060.9
$ws.Range("N68").Value = -3498
$ws.Range("H71").Value = 1818.9524
$ws.Range("I71").Value = 1809.9
$ws.Range("J71").Value = 2000
$ws.Range("K71").Value = 9049.5
$ws.Range("L71").Value = 10000
$ws.Range("M71").Value = -5305.5
$ws.Range("N71").Value = -17488
